$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "Nome:" / NOME field placeholder
$findRng = $d.Content
$found = $findRng.Find.Execute("Nome:")
if (-not $found) {
    throw "Could not find target paragraph ('Nome:') to replace"
}
$targetPara = $findRng.Paragraphs(1)
$targetRange = $targetPara.Range

# Build the replacement paragraph: a free-text lead-in ("Eu,") followed by four
# Word text form fields (NOME, SOBRENOME, CPF, RG), each one a bookmarked
# FORMTEXT field with its own w:ffData, mirroring the "use-the-template" wiring
# added in this revision.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Eu,</w:t></w:r><w:r><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val="NOME"/><w:enabled/><w:calcOnExit w:val="0"/><w:helpText w:type="autoText" w:val="NOME"/><w:textInput><w:default w:val="NOME"/></w:textInput></w:ffData></w:fldChar></w:r><w:bookmarkStart w:id="0" w:name="NOME"/><w:r><w:instrText xml:space="preserve"> FORMTEXT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:t>[ NOME ]</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val="SOBRENOME"/><w:enabled/><w:calcOnExit w:val="0"/><w:helpText w:type="autoText" w:val="SOBRENOME"/><w:textInput><w:default w:val="SOBRENOME"/></w:textInput></w:ffData></w:fldChar></w:r><w:bookmarkStart w:id="1" w:name="SOBRENOME"/><w:r><w:instrText xml:space="preserve"> FORMTEXT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:t>[ SOBRENOME ]</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:bookmarkEnd w:id="1"/><w:r><w:t xml:space="preserve">, com número de CPF </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val="CPF"/><w:enabled/><w:calcOnExit w:val="0"/><w:helpText w:type="autoText" w:val="CPF"/><w:textInput><w:default w:val="CPF"/></w:textInput></w:ffData></w:fldChar></w:r><w:bookmarkStart w:id="2" w:name="CPF"/><w:r><w:instrText xml:space="preserve"> FORMTEXT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:t>[ CPF ]</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:bookmarkEnd w:id="2"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">e RG </w:t></w:r><w:r><w:fldChar w:fldCharType="begin"><w:ffData><w:name w:val="RG"/><w:enabled/><w:calcOnExit w:val="0"/><w:helpText w:type="autoText" w:val="RG"/><w:textInput><w:default w:val="RG"/></w:textInput></w:ffData></w:fldChar></w:r><w:bookmarkStart w:id="3" w:name="RG"/><w:r><w:instrText xml:space="preserve"> FORMTEXT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:t>[ RG ]</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:bookmarkStart w:id="4" w:name="_GoBack"/><w:bookmarkEnd w:id="3"/><w:bookmarkEnd w:id="4"/></w:p>'

# InsertXML replaces the full contents of the addressed range with the supplied
# OOXML fragment, so this swaps the old "Nome: { NOME }" field-code paragraph for
# the new one in a single shot.
$targetRange.InsertXML($xml)

Write-Output "Paragraph replaced with NOME/SOBRENOME/CPF/RG form fields"
